$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: backfill column A (id) for existing rows 353..417 ---
# Pattern observed elsewhere in the sheet: A<row> = <row> - 2
for ($r = 353; $r -le 417; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Step 2: append the new "Summer 2016" rows (418..441) ---
# Columns: A=id, B=title, C=status, D=priority, E=season, F=year, G=remarks
$newRows = @(
    @(416, "Danganronpa 3 - The End of Kibougamine Gakuen - Mirai-hen", 1, 3, 2, 2016, $null),
    @(417, "Danganronpa 3 - The End of Kibougamine Gakuen - Zetsubou-hen", 1, 3, 2, 2016, $null),
    @(418, "Kono Bijutsu-bu ni wa Mondai ga Aru!", 1, 3, 2, 2016, $null),
    @(419, "Nanatsu no Taizai - Seisen no Shirushi", 2, 3, 2, 2016, $null),
    @(420, "ReLIFE", 1, 3, 2, 2016, "WAITING FOR FFF BDs (Partially Downloaded)"),
    @(421, "Time Travel Shoujo - Mari Waka to 8-nin no Kagakusha-tachi", 1, 3, 2, 2016, $null),
    @(422, "Ange Vierge", 1, 2, 2, 2016, $null),
    @(423, "Hatsukoi Monster", 1, 2, 2, 2016, $null),
    @(424, "New Game!", 1, 2, 2, 2016, $null),
    @(425, "Rewrite", 1, 2, 2, 2016, $null),
    @(426, "Fukigen na Mononokean", 1, 2, 2, 2016, $null),
    @(427, "Masou Gakuen HxH", 1, 2, 2, 2016, $null),
    @(428, "Momokuri", 1, 2, 2, 2016, $null),
    @(429, "Nejimaki Seirei Senki - Tenkyou no Alderamin", 1, 2, 2, 2016, $null),
    @(430, "Taboo Tattoo", 1, 2, 2, 2016, $null),
    @(431, "91Days", 1, 1, 2, 2016, $null),
    @(432, "Amaama to Inazuma", 1, 1, 2, 2016, $null),
    @(433, "Amanchu!", 1, 1, 2, 2016, $null),
    @(434, "Fate/kaleid liner Prisma Illya 3rei!!", 1, 1, 2, 2016, $null),
    @(435, "Hitori no Shita - The Outcast", 1, 1, 2, 2016, $null),
    @(436, "Orange", 3, 1, 2, 2016, $null),
    @(437, "Regalia - The Three Sacred Stars", 1, 1, 2, 2016, $null),
    @(438, "Scared Rider XechS", 1, 1, 2, 2016, $null),
    @(439, "Servamp", 1, 1, 2, 2016, $null)
)

# First pass: write id/title/status/priority/season/year for every new row so
# the "title" shared strings land contiguously (matching the authoring order).
$r = 418
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Second pass: fill in any "remarks" (column G) only after all titles exist,
# so the new shared string is appended last (after the 24 new titles).
$r = 418
foreach ($row in $newRows) {
    if ($row[6] -ne $null) {
        $ws.Cells.Item($r, 7).Value = $row[6]
    }
    $r++
}

# --- Step 3: restore the frozen-pane view / selection to match the new extent ---
$ws.Range("A391").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B416").Select()
